# BehaviorScenario_Technology_Power.xlsx
# "annual electricity consumption reduced to 1200+"
#
# Two data rows (ID_Technology = 22 and 23, the original rows 23 & 24) were
# removed from the table; the remaining rows shift up and the ID_Technology
# column is renumbered sequentially again (1..37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the two rows (this also shifts every cell/style below them up,
#    which is exactly what happened to the "unit"/"value" columns).
$ws.Rows("23:24").Delete()

# 2) Renumber the ID_Technology column (column A) back to a clean 1..37
#    sequence for the remaining data rows (2..38).
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# 3) Column A carries its own banded-fill styling that doesn't travel with
#    the deleted rows the way column C's does - restore the alternating
#    shaded/unshaded pattern for rows 23..38 to match the final layout.
#    Rows 35..38 also used to carry the "last rows" orange highlight font
#    (it shifted up from the old rows 37..40); column A for those rows goes
#    back to the regular black font now that row 38 is the true last row.
$addShadeRowsA  = @(25, 27, 31, 33, 35)
$removeShadeRowsA = @(36, 38)
$blackFontRowsA  = @(35, 36, 37, 38)

foreach ($r in $addShadeRowsA) {
    $ws.Cells.Item($r, 1).Interior.Color = 15132391
}
foreach ($r in $removeShadeRowsA) {
    $ws.Cells.Item($r, 1).Interior.ColorIndex = -4142
    $ws.Cells.Item($r, 1).Interior.Pattern = -4142
}
foreach ($r in $blackFontRowsA) {
    $ws.Cells.Item($r, 1).Font.Color = 0
}

# 4) Sheet view cosmetics: no forced top-left scroll anchor anymore and the
#    active selection lands on C32.
$ws.Range("A1").Select()
$ws.Range("C32").Select()
